# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market data refresh to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 322.06668
$ws.Range("I4").Value = 275.07144
$ws.Range("J4").Value = 980
$ws.Range("K4").Value = 275.07144
$ws.Range("L4").Value = 980
$ws.Range("M4").Value = -161.07144
$ws.Range("N4").Value = -1208
$ws.Range("H103").Value = 531.8889
$ws.Range("J103").Value = 531.3333
$ws.Range("L103").Value = 1593.9999
$ws.Range("N103").Value = -2765.9999
$ws.Range("H113").Value = 2434.3333
$ws.Range("J113").Value = 2427.8333
$ws.Range("L113").Value = 2427.8333
$ws.Range("N113").Value = -8935.8333
$ws.Range("H127").Value = 992.61536
$ws.Range("I127").Value = 547.7895
$ws.Range("J127").Value = 2200
$ws.Range("K127").Value = 1643.3685
$ws.Range("L127").Value = 6600
$ws.Range("M127").Value = 3316.6315
$ws.Range("N127").Value = -16520
$ws.Range("H132").Value = 6066908.5
$ws.Range("I132").Value = 10421833
$ws.Range("J132").Value = 7883.174
$ws.Range("K132").Value = 31265499
$ws.Range("L132").Value = 23649.522
$ws.Range("M132").Value = -31262969
$ws.Range("N132").Value = -28709.522
$ws.Range("H135").Value = 598.27026
$ws.Range("I135").Value = 185.93939
$ws.Range("K135").Value = 1673.45451
$ws.Range("M135").Value = 861.54549
$ws.Range("H137").Value = 1050.735
$ws.Range("I137").Value = 861.93475
$ws.Range("J137").Value = 1285.4595
$ws.Range("K137").Value = 2585.80425
$ws.Range("L137").Value = 3856.3785
$ws.Range("M137").Value = -35.80425000000014
$ws.Range("N137").Value = -8956.378499999999
$ws.Range("H138").Value = 622181.5600000001
$ws.Range("I138").Value = 723.54346
$ws.Range("J138").Value = 1813309.4
$ws.Range("K138").Value = 2170.63038
$ws.Range("L138").Value = 5439928.199999999
$ws.Range("M138").Value = 2969.36962
$ws.Range("N138").Value = -5450208.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 23256764
$ws.Range("I61").Value = 28572316
$ws.Range("J61").Value = 1226.875
$ws.Range("K61").Value = 28572316
$ws.Range("L61").Value = 1226.875
$ws.Range("M61").Value = -28572104
$ws.Range("N61").Value = -1650.875
$ws.Range("H74").Value = 789.65
$ws.Range("I74").Value = 672.9474
$ws.Range("K74").Value = 672.9474
$ws.Range("M74").Value = 201.0526
$ws.Range("H77").Value = 789.65
$ws.Range("I77").Value = 672.9474
$ws.Range("K77").Value = 3364.737
$ws.Range("M77").Value = 1003.263
$ws.Range("H132").Value = 2180.5652
$ws.Range("I132").Value = 2158.75
$ws.Range("J132").Value = 2230.4285
$ws.Range("K132").Value = 6476.25
$ws.Range("L132").Value = 6691.2855
$ws.Range("M132").Value = -3946.25
$ws.Range("N132").Value = -11751.2855
$ws.Range("H136").Value = 23256764
$ws.Range("I136").Value = 28572316
$ws.Range("J136").Value = 1226.875
$ws.Range("K136").Value = 85716948
$ws.Range("L136").Value = 3680.625
$ws.Range("M136").Value = -85714398
$ws.Range("N136").Value = -8780.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1250.5
$ws.Range("I22").Value = 1250.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1250.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1077.5
$ws.Range("N22").Value = ""
$ws.Range("H134").Value = 4630.385
$ws.Range("I134").Value = 1366.2122
$ws.Range("J134").Value = 22583.334
$ws.Range("K134").Value = 4098.6366
$ws.Range("L134").Value = 67750.00199999999
$ws.Range("M134").Value = -1563.6366
$ws.Range("N134").Value = -72820.00199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 321
$ws.Range("I7").Value = 191
$ws.Range("J7").Value = 477
$ws.Range("K7").Value = 191
$ws.Range("L7").Value = 477
$ws.Range("M7").Value = -78
$ws.Range("N7").Value = -703
$ws.Range("H16").Value = 83334856
$ws.Range("J16").Value = 1412.5
$ws.Range("L16").Value = 1412.5
$ws.Range("N16").Value = -1986.5
$ws.Range("H31").Value = 1797.3438
$ws.Range("I31").Value = 1925.0385
$ws.Range("K31").Value = 1925.0385
$ws.Range("M31").Value = -1630.0385
$ws.Range("H34").Value = 1797.3438
$ws.Range("I34").Value = 1925.0385
$ws.Range("K34").Value = 1925.0385
$ws.Range("M34").Value = -1723.0385
$ws.Range("H113").Value = 83334856
$ws.Range("J113").Value = 1412.5
$ws.Range("L113").Value = 1412.5
$ws.Range("N113").Value = -5752.5
$ws.Range("H114").Value = 26245
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 26245
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 26245
$ws.Range("M114").Value = ""
$ws.Range("N114").Value = -34923
$ws.Range("H122").Value = 1151.875
$ws.Range("I122").Value = 992.36365
$ws.Range("K122").Value = 2977.09095
$ws.Range("M122").Value = -527.0909499999998
$ws.Range("H132").Value = 2985.635
$ws.Range("I132").Value = 2888.3333
$ws.Range("J132").Value = 3569.4443
$ws.Range("K132").Value = 8664.999899999999
$ws.Range("L132").Value = 10708.3329
$ws.Range("M132").Value = -6134.999899999999
$ws.Range("N132").Value = -15768.3329
$ws.Range("H134").Value = 9091982
$ws.Range("I134").Value = 1093.6459
$ws.Range("J134").Value = 71429500
$ws.Range("K134").Value = 3280.9377
$ws.Range("L134").Value = 214288500
$ws.Range("M134").Value = -745.9376999999999
$ws.Range("N134").Value = -214293570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1744.25
$ws.Range("I5").Value = 1920.1111
$ws.Range("J5").Value = 1216.6666
$ws.Range("K5").Value = 5760.3333
$ws.Range("L5").Value = 3649.9998
$ws.Range("M5").Value = -5648.3333
$ws.Range("N5").Value = -3873.9998
$ws.Range("H34").Value = 1797.625
$ws.Range("J34").Value = 2055.8333
$ws.Range("L34").Value = 6167.499899999999
$ws.Range("N34").Value = -6335.499899999999
$ws.Range("H80").Value = 3690
$ws.Range("H83").Value = 3690
$ws.Range("H121").Value = 525.4286
$ws.Range("J121").Value = 682
$ws.Range("L121").Value = 2046
$ws.Range("N121").Value = -4666
$ws.Range("H131").Value = 16667890
$ws.Range("J131").Value = 1328.5471
$ws.Range("L131").Value = 3985.6413
$ws.Range("N131").Value = -14065.6413
$ws.Range("H132").Value = 2040.8
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060
$ws.Range("H135").Value = 1744.25
$ws.Range("I135").Value = 1920.1111
$ws.Range("J135").Value = 1216.6666
$ws.Range("K135").Value = 17280.9999
$ws.Range("L135").Value = 10949.9994
$ws.Range("M135").Value = -14745.9999
$ws.Range("N135").Value = -16019.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 558.7778
$ws.Range("I97").Value = 558.7778
$ws.Range("K97").Value = 558.7778
$ws.Range("M97").Value = -62.77779999999996
$ws.Range("H113").Value = 1235.1666
$ws.Range("I113").Value = 1342.2
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 1342.2
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 827.8
$ws.Range("N113").Value = -5040
$ws.Range("H122").Value = 1860.8
$ws.Range("I122").Value = 1446.6364
$ws.Range("J122").Value = 2999.75
$ws.Range("K122").Value = 4339.9092
$ws.Range("L122").Value = 8999.25
$ws.Range("M122").Value = -1889.9092
$ws.Range("N122").Value = -13899.25
$ws.Range("H132").Value = 1996.2812
$ws.Range("I132").Value = 1742.174
$ws.Range("K132").Value = 5226.522
$ws.Range("M132").Value = -2696.522

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
$ws.Range("H122").Value = 31268876
$ws.Range("I122").Value = 62525900
$ws.Range("K122").Value = 187577700
$ws.Range("M122").Value = -187575250
$ws.Range("H132").Value = 20836.36
$ws.Range("I132").Value = 1540.4667
$ws.Range("J132").Value = 46004.914
$ws.Range("K132").Value = 4621.4001
$ws.Range("L132").Value = 138014.742
$ws.Range("M132").Value = -2091.4001
$ws.Range("N132").Value = -143074.742
$ws.Range("H136").Value = 2999.681
$ws.Range("I136").Value = 2999.689
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 8999.066999999999
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -6449.066999999999
$ws.Range("N136").Value = -14098.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3000
$ws.Range("J4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("N4").Value = -3226
$ws.Range("H100").Value = 705.1667
$ws.Range("J100").Value = 647
$ws.Range("L100").Value = 1294
$ws.Range("N100").Value = -2376
$ws.Range("H107").Value = 444.72223
$ws.Range("I107").Value = 447.91666
$ws.Range("K107").Value = 1343.74998
$ws.Range("M107").Value = 576.2500199999999
$ws.Range("H113").Value = 315.95456
$ws.Range("I113").Value = 208
$ws.Range("K113").Value = 624
$ws.Range("M113").Value = 1546
$ws.Range("H122").Value = 89288130
$ws.Range("I122").Value = 113638350
$ws.Range("K122").Value = 340915050
$ws.Range("M122").Value = -340912600
$ws.Range("H126").Value = 30304258
$ws.Range("I126").Value = 40001104
$ws.Range("J126").Value = 1612.5
$ws.Range("K126").Value = 120003312
$ws.Range("L126").Value = 4837.5
$ws.Range("M126").Value = -120000842
$ws.Range("N126").Value = -9777.5
$ws.Range("H132").Value = 1813.1724
$ws.Range("I132").Value = 1758.9807
$ws.Range("K132").Value = 5276.9421
$ws.Range("M132").Value = -2746.9421
$ws.Range("H136").Value = 528.34045
$ws.Range("I136").Value = 406.82352
$ws.Range("J136").Value = 846.1539
$ws.Range("K136").Value = 1220.47056
$ws.Range("L136").Value = 2538.4617
$ws.Range("M136").Value = 1329.52944
$ws.Range("N136").Value = -7638.4617
